# Update the "AYKO" interactive-map workbook:
#  - Remove the data row for Caso 6398 (LARRAZABAL AV. 579) - originally row 54
#  - Remove the data row for Caso -617 (Soler 3815) - originally row 78
#  - Remove the now-unused "PD" (Q) and "N2" (R) columns entirely
#
# Deleting rows/columns shifts everything below/right of them up/left,
# so the sheet's used range goes from A1:R85 to A1:P83.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the higher-numbered row first so the row 54 index stays valid.
$ws.Rows.Item(78).Delete()
$ws.Rows.Item(54).Delete()

# Columns Q ("PD") and R ("N2") are removed completely.
$ws.Range("Q:R").Delete()
